# Update metrics_2_5 results: re-rank models (4 melhores rstds) by
# rewriting model names in column A (new row order) and refreshing the
# metric columns B:I for every data row (rows 2-26) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_2_5_0"
$ws.Cells.Item(2, 2).Value = 0.09932080507882668
$ws.Cells.Item(2, 3).Value = -0.8325171248060477
$ws.Cells.Item(2, 4).Value = 0.04054228265915827
$ws.Cells.Item(2, 5).Value = -0.006708787800333571
$ws.Cells.Item(2, 6).Value = 0.9967864155769348
$ws.Cells.Item(2, 7).Value = 1.248383045196533
$ws.Cells.Item(2, 8).Value = 1.570474863052368
$ws.Cells.Item(2, 9).Value = 1.399955749511719
$ws.Cells.Item(3, 1).Value = "model_2_5_22"
$ws.Cells.Item(3, 2).Value = 0.09932080507882668
$ws.Cells.Item(3, 3).Value = -0.8325171248060477
$ws.Cells.Item(3, 4).Value = 0.04054228265915827
$ws.Cells.Item(3, 5).Value = -0.006708787800333571
$ws.Cells.Item(3, 6).Value = 0.9967864155769348
$ws.Cells.Item(3, 7).Value = 1.248383045196533
$ws.Cells.Item(3, 8).Value = 1.570474863052368
$ws.Cells.Item(3, 9).Value = 1.399955749511719
$ws.Cells.Item(4, 1).Value = "model_2_5_21"
$ws.Cells.Item(4, 2).Value = 0.09932080507882668
$ws.Cells.Item(4, 3).Value = -0.8325171248060477
$ws.Cells.Item(4, 4).Value = 0.04054228265915827
$ws.Cells.Item(4, 5).Value = -0.006708787800333571
$ws.Cells.Item(4, 6).Value = 0.9967864155769348
$ws.Cells.Item(4, 7).Value = 1.248383045196533
$ws.Cells.Item(4, 8).Value = 1.570474863052368
$ws.Cells.Item(4, 9).Value = 1.399955749511719
$ws.Cells.Item(5, 1).Value = "model_2_5_20"
$ws.Cells.Item(5, 2).Value = 0.09932080507882668
$ws.Cells.Item(5, 3).Value = -0.8325171248060477
$ws.Cells.Item(5, 4).Value = 0.04054228265915827
$ws.Cells.Item(5, 5).Value = -0.006708787800333571
$ws.Cells.Item(5, 6).Value = 0.9967864155769348
$ws.Cells.Item(5, 7).Value = 1.248383045196533
$ws.Cells.Item(5, 8).Value = 1.570474863052368
$ws.Cells.Item(5, 9).Value = 1.399955749511719
$ws.Cells.Item(6, 1).Value = "model_2_5_19"
$ws.Cells.Item(6, 2).Value = 0.09932080507882668
$ws.Cells.Item(6, 3).Value = -0.8325171248060477
$ws.Cells.Item(6, 4).Value = 0.04054228265915827
$ws.Cells.Item(6, 5).Value = -0.006708787800333571
$ws.Cells.Item(6, 6).Value = 0.9967864155769348
$ws.Cells.Item(6, 7).Value = 1.248383045196533
$ws.Cells.Item(6, 8).Value = 1.570474863052368
$ws.Cells.Item(6, 9).Value = 1.399955749511719
$ws.Cells.Item(7, 1).Value = "model_2_5_18"
$ws.Cells.Item(7, 2).Value = 0.09932080507882668
$ws.Cells.Item(7, 3).Value = -0.8325171248060477
$ws.Cells.Item(7, 4).Value = 0.04054228265915827
$ws.Cells.Item(7, 5).Value = -0.006708787800333571
$ws.Cells.Item(7, 6).Value = 0.9967864155769348
$ws.Cells.Item(7, 7).Value = 1.248383045196533
$ws.Cells.Item(7, 8).Value = 1.570474863052368
$ws.Cells.Item(7, 9).Value = 1.399955749511719
$ws.Cells.Item(8, 1).Value = "model_2_5_17"
$ws.Cells.Item(8, 2).Value = 0.09932080507882668
$ws.Cells.Item(8, 3).Value = -0.8325171248060477
$ws.Cells.Item(8, 4).Value = 0.04054228265915827
$ws.Cells.Item(8, 5).Value = -0.006708787800333571
$ws.Cells.Item(8, 6).Value = 0.9967864155769348
$ws.Cells.Item(8, 7).Value = 1.248383045196533
$ws.Cells.Item(8, 8).Value = 1.570474863052368
$ws.Cells.Item(8, 9).Value = 1.399955749511719
$ws.Cells.Item(9, 1).Value = "model_2_5_16"
$ws.Cells.Item(9, 2).Value = 0.09932080507882668
$ws.Cells.Item(9, 3).Value = -0.8325171248060477
$ws.Cells.Item(9, 4).Value = 0.04054228265915827
$ws.Cells.Item(9, 5).Value = -0.006708787800333571
$ws.Cells.Item(9, 6).Value = 0.9967864155769348
$ws.Cells.Item(9, 7).Value = 1.248383045196533
$ws.Cells.Item(9, 8).Value = 1.570474863052368
$ws.Cells.Item(9, 9).Value = 1.399955749511719
$ws.Cells.Item(10, 1).Value = "model_2_5_15"
$ws.Cells.Item(10, 2).Value = 0.09932080507882668
$ws.Cells.Item(10, 3).Value = -0.8325171248060477
$ws.Cells.Item(10, 4).Value = 0.04054228265915827
$ws.Cells.Item(10, 5).Value = -0.006708787800333571
$ws.Cells.Item(10, 6).Value = 0.9967864155769348
$ws.Cells.Item(10, 7).Value = 1.248383045196533
$ws.Cells.Item(10, 8).Value = 1.570474863052368
$ws.Cells.Item(10, 9).Value = 1.399955749511719
$ws.Cells.Item(11, 1).Value = "model_2_5_14"
$ws.Cells.Item(11, 2).Value = 0.09932080507882668
$ws.Cells.Item(11, 3).Value = -0.8325171248060477
$ws.Cells.Item(11, 4).Value = 0.04054228265915827
$ws.Cells.Item(11, 5).Value = -0.006708787800333571
$ws.Cells.Item(11, 6).Value = 0.9967864155769348
$ws.Cells.Item(11, 7).Value = 1.248383045196533
$ws.Cells.Item(11, 8).Value = 1.570474863052368
$ws.Cells.Item(11, 9).Value = 1.399955749511719
$ws.Cells.Item(12, 1).Value = "model_2_5_13"
$ws.Cells.Item(12, 2).Value = 0.09932080507882668
$ws.Cells.Item(12, 3).Value = -0.8325171248060477
$ws.Cells.Item(12, 4).Value = 0.04054228265915827
$ws.Cells.Item(12, 5).Value = -0.006708787800333571
$ws.Cells.Item(12, 6).Value = 0.9967864155769348
$ws.Cells.Item(12, 7).Value = 1.248383045196533
$ws.Cells.Item(12, 8).Value = 1.570474863052368
$ws.Cells.Item(12, 9).Value = 1.399955749511719
$ws.Cells.Item(13, 1).Value = "model_2_5_23"
$ws.Cells.Item(13, 2).Value = 0.09932080507882668
$ws.Cells.Item(13, 3).Value = -0.8325171248060477
$ws.Cells.Item(13, 4).Value = 0.04054228265915827
$ws.Cells.Item(13, 5).Value = -0.006708787800333571
$ws.Cells.Item(13, 6).Value = 0.9967864155769348
$ws.Cells.Item(13, 7).Value = 1.248383045196533
$ws.Cells.Item(13, 8).Value = 1.570474863052368
$ws.Cells.Item(13, 9).Value = 1.399955749511719
$ws.Cells.Item(14, 1).Value = "model_2_5_12"
$ws.Cells.Item(14, 2).Value = 0.09932080507882668
$ws.Cells.Item(14, 3).Value = -0.8325171248060477
$ws.Cells.Item(14, 4).Value = 0.04054228265915827
$ws.Cells.Item(14, 5).Value = -0.006708787800333571
$ws.Cells.Item(14, 6).Value = 0.9967864155769348
$ws.Cells.Item(14, 7).Value = 1.248383045196533
$ws.Cells.Item(14, 8).Value = 1.570474863052368
$ws.Cells.Item(14, 9).Value = 1.399955749511719
$ws.Cells.Item(15, 1).Value = "model_2_5_10"
$ws.Cells.Item(15, 2).Value = 0.09932080507882668
$ws.Cells.Item(15, 3).Value = -0.8325171248060477
$ws.Cells.Item(15, 4).Value = 0.04054228265915827
$ws.Cells.Item(15, 5).Value = -0.006708787800333571
$ws.Cells.Item(15, 6).Value = 0.9967864155769348
$ws.Cells.Item(15, 7).Value = 1.248383045196533
$ws.Cells.Item(15, 8).Value = 1.570474863052368
$ws.Cells.Item(15, 9).Value = 1.399955749511719
$ws.Cells.Item(16, 1).Value = "model_2_5_9"
$ws.Cells.Item(16, 2).Value = 0.09932080507882668
$ws.Cells.Item(16, 3).Value = -0.8325171248060477
$ws.Cells.Item(16, 4).Value = 0.04054228265915827
$ws.Cells.Item(16, 5).Value = -0.006708787800333571
$ws.Cells.Item(16, 6).Value = 0.9967864155769348
$ws.Cells.Item(16, 7).Value = 1.248383045196533
$ws.Cells.Item(16, 8).Value = 1.570474863052368
$ws.Cells.Item(16, 9).Value = 1.399955749511719
$ws.Cells.Item(17, 1).Value = "model_2_5_8"
$ws.Cells.Item(17, 2).Value = 0.09932080507882668
$ws.Cells.Item(17, 3).Value = -0.8325171248060477
$ws.Cells.Item(17, 4).Value = 0.04054228265915827
$ws.Cells.Item(17, 5).Value = -0.006708787800333571
$ws.Cells.Item(17, 6).Value = 0.9967864155769348
$ws.Cells.Item(17, 7).Value = 1.248383045196533
$ws.Cells.Item(17, 8).Value = 1.570474863052368
$ws.Cells.Item(17, 9).Value = 1.399955749511719
$ws.Cells.Item(18, 1).Value = "model_2_5_7"
$ws.Cells.Item(18, 2).Value = 0.09932080507882668
$ws.Cells.Item(18, 3).Value = -0.8325171248060477
$ws.Cells.Item(18, 4).Value = 0.04054228265915827
$ws.Cells.Item(18, 5).Value = -0.006708787800333571
$ws.Cells.Item(18, 6).Value = 0.9967864155769348
$ws.Cells.Item(18, 7).Value = 1.248383045196533
$ws.Cells.Item(18, 8).Value = 1.570474863052368
$ws.Cells.Item(18, 9).Value = 1.399955749511719
$ws.Cells.Item(19, 1).Value = "model_2_5_6"
$ws.Cells.Item(19, 2).Value = 0.09932080507882668
$ws.Cells.Item(19, 3).Value = -0.8325171248060477
$ws.Cells.Item(19, 4).Value = 0.04054228265915827
$ws.Cells.Item(19, 5).Value = -0.006708787800333571
$ws.Cells.Item(19, 6).Value = 0.9967864155769348
$ws.Cells.Item(19, 7).Value = 1.248383045196533
$ws.Cells.Item(19, 8).Value = 1.570474863052368
$ws.Cells.Item(19, 9).Value = 1.399955749511719
$ws.Cells.Item(20, 1).Value = "model_2_5_5"
$ws.Cells.Item(20, 2).Value = 0.09932080507882668
$ws.Cells.Item(20, 3).Value = -0.8325171248060477
$ws.Cells.Item(20, 4).Value = 0.04054228265915827
$ws.Cells.Item(20, 5).Value = -0.006708787800333571
$ws.Cells.Item(20, 6).Value = 0.9967864155769348
$ws.Cells.Item(20, 7).Value = 1.248383045196533
$ws.Cells.Item(20, 8).Value = 1.570474863052368
$ws.Cells.Item(20, 9).Value = 1.399955749511719
$ws.Cells.Item(21, 1).Value = "model_2_5_4"
$ws.Cells.Item(21, 2).Value = 0.09932080507882668
$ws.Cells.Item(21, 3).Value = -0.8325171248060477
$ws.Cells.Item(21, 4).Value = 0.04054228265915827
$ws.Cells.Item(21, 5).Value = -0.006708787800333571
$ws.Cells.Item(21, 6).Value = 0.9967864155769348
$ws.Cells.Item(21, 7).Value = 1.248383045196533
$ws.Cells.Item(21, 8).Value = 1.570474863052368
$ws.Cells.Item(21, 9).Value = 1.399955749511719
$ws.Cells.Item(22, 1).Value = "model_2_5_3"
$ws.Cells.Item(22, 2).Value = 0.09932080507882668
$ws.Cells.Item(22, 3).Value = -0.8325171248060477
$ws.Cells.Item(22, 4).Value = 0.04054228265915827
$ws.Cells.Item(22, 5).Value = -0.006708787800333571
$ws.Cells.Item(22, 6).Value = 0.9967864155769348
$ws.Cells.Item(22, 7).Value = 1.248383045196533
$ws.Cells.Item(22, 8).Value = 1.570474863052368
$ws.Cells.Item(22, 9).Value = 1.399955749511719
$ws.Cells.Item(23, 1).Value = "model_2_5_2"
$ws.Cells.Item(23, 2).Value = 0.09932080507882668
$ws.Cells.Item(23, 3).Value = -0.8325171248060477
$ws.Cells.Item(23, 4).Value = 0.04054228265915827
$ws.Cells.Item(23, 5).Value = -0.006708787800333571
$ws.Cells.Item(23, 6).Value = 0.9967864155769348
$ws.Cells.Item(23, 7).Value = 1.248383045196533
$ws.Cells.Item(23, 8).Value = 1.570474863052368
$ws.Cells.Item(23, 9).Value = 1.399955749511719
$ws.Cells.Item(24, 1).Value = "model_2_5_1"
$ws.Cells.Item(24, 2).Value = 0.09932080507882668
$ws.Cells.Item(24, 3).Value = -0.8325171248060477
$ws.Cells.Item(24, 4).Value = 0.04054228265915827
$ws.Cells.Item(24, 5).Value = -0.006708787800333571
$ws.Cells.Item(24, 6).Value = 0.9967864155769348
$ws.Cells.Item(24, 7).Value = 1.248383045196533
$ws.Cells.Item(24, 8).Value = 1.570474863052368
$ws.Cells.Item(24, 9).Value = 1.399955749511719
$ws.Cells.Item(25, 1).Value = "model_2_5_11"
$ws.Cells.Item(25, 2).Value = 0.09932080507882668
$ws.Cells.Item(25, 3).Value = -0.8325171248060477
$ws.Cells.Item(25, 4).Value = 0.04054228265915827
$ws.Cells.Item(25, 5).Value = -0.006708787800333571
$ws.Cells.Item(25, 6).Value = 0.9967864155769348
$ws.Cells.Item(25, 7).Value = 1.248383045196533
$ws.Cells.Item(25, 8).Value = 1.570474863052368
$ws.Cells.Item(25, 9).Value = 1.399955749511719
$ws.Cells.Item(26, 1).Value = "model_2_5_24"
$ws.Cells.Item(26, 2).Value = 0.09932080507882668
$ws.Cells.Item(26, 3).Value = -0.8325171248060477
$ws.Cells.Item(26, 4).Value = 0.04054228265915827
$ws.Cells.Item(26, 5).Value = -0.006708787800333571
$ws.Cells.Item(26, 6).Value = 0.9967864155769348
$ws.Cells.Item(26, 7).Value = 1.248383045196533
$ws.Cells.Item(26, 8).Value = 1.570474863052368
$ws.Cells.Item(26, 9).Value = 1.399955749511719
